# Add a new check-list test case (GW_3: "Check that the button for
# search is working") and extend the Expected Results text of GW_2 / GW_3
# with an extra verification step about the entered text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Row 3 (GW_2 test case): update "Test Case Description" (E3) to add
#    a 3rd numbered step, and grow the row height to fit the new text.
# ------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 76.5

$e3Bold = "Check that the search is working"
$e3Rest = "`n1. Open 'https://www.google.com/' with Chrome`n2. Enter a valid query in the search (for example: wikipedia). `n3. Use key 'Enter'. "
$ws.Range("E3").Value = $e3Bold + $e3Rest
$ws.Range("E3").Characters(1, $e3Bold.Length).Font.Bold = $true
$e3RestChars = $ws.Range("E3").Characters($e3Bold.Length + 1, $e3Rest.Length)
$e3RestChars.Font.Bold = $false
$e3RestChars.Font.Name = "Calibri"
$e3RestChars.Font.Size = 11

# ------------------------------------------------------------------
# 2. Row 3, column F ("Expected Results") used to hold the long expected
#    -results text for GW_2; that text now becomes shared with the new
#    GW_3 row below, while F3 itself is re-pointed at it (same string
#    the new F4 will use) after the extra "entered text" verification
#    step was appended.
# ------------------------------------------------------------------
$expectedResults = "`n1. Website correctly open on Chrome browser.`n2. The entered text is displayed correctly in the input field.`n3. Google search page with query results correctly open."
$ws.Range("F3").Value = $expectedResults

# ------------------------------------------------------------------
# 3. Row 4 (new GW_3 test case) - copy formatting from row 3's cells so
#    styles match exactly, then fill in the new content.
# ------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 81.75

$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A4").Value = "GW_3"

$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)   # xlPasteFormats
$e4Bold = "Check that the button for search is working"
$e4Rest = "`n1. Open 'https://www.google.com/' with Chrome`n2. Enter a valid query in the search (for example: wikipedia).`n3. Use button 'Google Search' for starting search . "
$ws.Range("E4").Value = $e4Bold + $e4Rest
$ws.Range("E4").Characters(1, $e4Bold.Length).Font.Bold = $true
$e4RestChars = $ws.Range("E4").Characters($e4Bold.Length + 1, $e4Rest.Length)
$e4RestChars.Font.Bold = $false
$e4RestChars.Font.Name = "Calibri"
$e4RestChars.Font.Size = 11

$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F4").Value = $expectedResults

$ws.Range("G4").Value = "Using english lang for browser, so that the name of the button displayed as in this test case description"

$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 4. Active selection moves from F4 to G4.
# ------------------------------------------------------------------
$ws.Range("G4").Select()
